# Auto-generated edit script: fixes mismatched match-result rows
# (team names / scores / odds swapped back onto their correct row)
# in the "Wales Premier League" sheet, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B98").Value = 6800493
$row98 = New-Object "object[,]" 1,24
$row98[0,0] = "Newtown"
$row98[0,1] = "Pontypridd Town"
$row98[0,2] = 3
$row98[0,3] = 1
$row98[0,4] = "H"
$row98[0,5] = 1.45
$row98[0,6] = 3.75
$row98[0,7] = 6.5
$row98[0,8] = 1.444
$row98[0,9] = 3.8
$row98[0,10] = 7.5
$row98[0,11] = -1.25
$row98[0,12] = 2.025
$row98[0,13] = 1.775
$row98[0,14] = 2.5
$row98[0,15] = 2
$row98[0,16] = 1.8
$row98[0,17] = 0.444
$row98[0,18] = -1
$row98[0,19] = -1
$row98[0,20] = 1.025
$row98[0,21] = -1
$row98[0,22] = 1
$row98[0,23] = -1
$ws.Range("F98:AC98").Value = $row98

$ws.Range("B99").Value = 6800495
$row99 = New-Object "object[,]" 1,24
$row99[0,0] = "Haverfordwest County"
$row99[0,1] = "Cardiff MU"
$row99[0,2] = 1
$row99[0,3] = 1
$row99[0,4] = "D"
$row99[0,5] = 2.5
$row99[0,6] = 3.4
$row99[0,7] = 2.4
$row99[0,8] = 2.45
$row99[0,9] = 3.2
$row99[0,10] = 2.7
$row99[0,11] = 0
$row99[0,12] = 1.775
$row99[0,13] = 2.025
$row99[0,14] = 2.5
$row99[0,15] = 1.95
$row99[0,16] = 1.85
$row99[0,17] = -1
$row99[0,18] = 2.2
$row99[0,19] = -1
$row99[0,20] = 0
$row99[0,21] = -0
$row99[0,22] = -1
$row99[0,23] = 0.8500000000000001
$ws.Range("F99:AC99").Value = $row99

$ws.Range("B100").Value = 6800492
$row100 = New-Object "object[,]" 1,24
$row100[0,0] = "Barry Town"
$row100[0,1] = "Penybont"
$row100[0,2] = 1
$row100[0,3] = 1
$row100[0,4] = "D"
$row100[0,5] = 3.75
$row100[0,6] = 3.75
$row100[0,7] = 1.727
$row100[0,8] = 3
$row100[0,9] = 3.6
$row100[0,10] = 2.05
$row100[0,11] = 0.25
$row100[0,12] = 2
$row100[0,13] = 1.8
$row100[0,14] = 2.75
$row100[0,15] = 1.875
$row100[0,16] = 1.925
$row100[0,17] = -1
$row100[0,18] = 2.6
$row100[0,19] = -1
$row100[0,20] = 0.5
$row100[0,21] = -0.5
$row100[0,22] = -1
$row100[0,23] = 0.925
$ws.Range("F100:AC100").Value = $row100

$ws.Range("B103").Value = 6800498
$row103 = New-Object "object[,]" 1,24
$row103[0,0] = "Newtown"
$row103[0,1] = "Colwyn Bay"
$row103[0,2] = 4
$row103[0,3] = 2
$row103[0,4] = "H"
$row103[0,5] = 1.444
$row103[0,6] = 4.2
$row103[0,7] = 6
$row103[0,8] = 1.4
$row103[0,9] = 4.2
$row103[0,10] = 6.5
$row103[0,11] = -1.25
$row103[0,12] = 1.925
$row103[0,13] = 1.875
$row103[0,14] = 3
$row103[0,15] = 1.925
$row103[0,16] = 1.875
$row103[0,17] = 0.3999999999999999
$row103[0,18] = -1
$row103[0,19] = -1
$row103[0,20] = 0.925
$row103[0,21] = -1
$row103[0,22] = 0.925
$row103[0,23] = -1
$ws.Range("F103:AC103").Value = $row103

$ws.Range("B104").Value = 6800497
$row104 = New-Object "object[,]" 1,24
$row104[0,0] = "Connahs Quay"
$row104[0,1] = "Barry Town"
$row104[0,2] = 7
$row104[0,3] = 0
$row104[0,4] = "H"
$row104[0,5] = 1.2
$row104[0,6] = 6.5
$row104[0,7] = 9
$row104[0,8] = 1.166
$row104[0,9] = 7.5
$row104[0,10] = 10
$row104[0,11] = -2
$row104[0,12] = 1.85
$row104[0,13] = 1.95
$row104[0,14] = 3.25
$row104[0,15] = 1.825
$row104[0,16] = 1.975
$row104[0,17] = 0.1659999999999999
$row104[0,18] = -1
$row104[0,19] = -1
$row104[0,20] = 0.8500000000000001
$row104[0,21] = -1
$row104[0,22] = 0.825
$row104[0,23] = -1
$ws.Range("F104:AC104").Value = $row104

$ws.Range("B108").Value = 6800049
$row108 = New-Object "object[,]" 1,24
$row108[0,0] = "Caernarfon Town"
$row108[0,1] = "Haverfordwest County"
$row108[0,2] = 0
$row108[0,3] = 1
$row108[0,4] = "A"
$row108[0,5] = 2.3
$row108[0,6] = 3.2
$row108[0,7] = 2.7
$row108[0,8] = 1.95
$row108[0,9] = 3.3
$row108[0,10] = 3.25
$row108[0,11] = -0.5
$row108[0,12] = 2.025
$row108[0,13] = 1.775
$row108[0,14] = 2.75
$row108[0,15] = 1.8
$row108[0,16] = 2
$row108[0,17] = -1
$row108[0,18] = -1
$row108[0,19] = 2.25
$row108[0,20] = -1
$row108[0,21] = 0.7749999999999999
$row108[0,22] = -1
$row108[0,23] = 1
$ws.Range("F108:AC108").Value = $row108

$ws.Range("B110").Value = 6800503
$row110 = New-Object "object[,]" 1,24
$row110[0,0] = "Connahs Quay"
$row110[0,1] = "Pontypridd Town"
$row110[0,2] = 3
$row110[0,3] = 1
$row110[0,4] = "H"
$row110[0,5] = 1.2
$row110[0,6] = 6
$row110[0,7] = 9
$row110[0,8] = 1.142
$row110[0,9] = 7
$row110[0,10] = 13
$row110[0,11] = -2
$row110[0,12] = 1.775
$row110[0,13] = 2.025
$row110[0,14] = 3.25
$row110[0,15] = 1.95
$row110[0,16] = 1.85
$row110[0,17] = 0.1419999999999999
$row110[0,18] = -1
$row110[0,19] = -1
$row110[0,20] = 0
$row110[0,21] = -0
$row110[0,22] = 0.95
$row110[0,23] = -1
$ws.Range("F110:AC110").Value = $row110

$ws.Range("B146").Value = 7721586
$row146 = New-Object "object[,]" 1,24
$row146[0,0] = "Caernarfon Town"
$row146[0,1] = "TNS"
$row146[0,2] = 1
$row146[0,3] = 8
$row146[0,4] = "A"
$row146[0,5] = 11
$row146[0,6] = 8
$row146[0,7] = 1.142
$row146[0,8] = 13
$row146[0,9] = 7.5
$row146[0,10] = 1.142
$row146[0,11] = 2.25
$row146[0,12] = 1.95
$row146[0,13] = 1.85
$row146[0,14] = 3.5
$row146[0,15] = 1.8
$row146[0,16] = 2
$row146[0,17] = -1
$row146[0,18] = -1
$row146[0,19] = 0.1419999999999999
$row146[0,20] = -1
$row146[0,21] = 0.8500000000000001
$row146[0,22] = 0.8
$row146[0,23] = -1
$ws.Range("F146:AC146").Value = $row146

$ws.Range("B147").Value = 7721608
$row147 = New-Object "object[,]" 1,24
$row147[0,0] = "Barry Town"
$row147[0,1] = "Haverfordwest County"
$row147[0,2] = 1
$row147[0,3] = 1
$row147[0,4] = "D"
$row147[0,5] = 2.3
$row147[0,6] = 3.4
$row147[0,7] = 2.75
$row147[0,8] = 2.25
$row147[0,9] = 3.25
$row147[0,10] = 2.9
$row147[0,11] = -0.25
$row147[0,12] = 2
$row147[0,13] = 1.8
$row147[0,14] = 2.5
$row147[0,15] = 1.925
$row147[0,16] = 1.875
$row147[0,17] = -1
$row147[0,18] = 2.25
$row147[0,19] = -1
$row147[0,20] = -0.5
$row147[0,21] = 0.4
$row147[0,22] = -1
$row147[0,23] = 0.875
$ws.Range("F147:AC147").Value = $row147

$ws.Range("B170").Value = 7721594
$row170 = New-Object "object[,]" 1,24
$row170[0,0] = "Connahs Quay"
$row170[0,1] = "Newtown"
$row170[0,2] = 0
$row170[0,3] = 0
$row170[0,4] = "D"
$row170[0,5] = 1.4
$row170[0,6] = 4.5
$row170[0,7] = 5.75
$row170[0,8] = 1.55
$row170[0,9] = 4.5
$row170[0,10] = 4.5
$row170[0,11] = -1
$row170[0,12] = 1.9
$row170[0,13] = 1.9
$row170[0,14] = 3
$row170[0,15] = 1.8
$row170[0,16] = 2
$row170[0,17] = -1
$row170[0,18] = 3.5
$row170[0,19] = -1
$row170[0,20] = -1
$row170[0,21] = 0.8999999999999999
$row170[0,22] = -1
$row170[0,23] = 1
$ws.Range("F170:AC170").Value = $row170

$ws.Range("B171").Value = 7721620
$row171 = New-Object "object[,]" 1,24
$row171[0,0] = "Barry Town"
$row171[0,1] = "Penybont"
$row171[0,2] = 0
$row171[0,3] = 0
$row171[0,4] = "D"
$row171[0,5] = 3.1
$row171[0,6] = 3.1
$row171[0,7] = 2.15
$row171[0,8] = 4
$row171[0,9] = 3.2
$row171[0,10] = 1.909
$row171[0,11] = 0.5
$row171[0,12] = 1.875
$row171[0,13] = 1.925
$row171[0,14] = 2.5
$row171[0,15] = 1.975
$row171[0,16] = 1.825
$row171[0,17] = -1
$row171[0,18] = 2.2
$row171[0,19] = -1
$row171[0,20] = 0.875
$row171[0,21] = -1
$row171[0,22] = -1
$row171[0,23] = 0.825
$ws.Range("F171:AC171").Value = $row171

$ws.Range("B172").Value = 7721621
$row172 = New-Object "object[,]" 1,24
$row172[0,0] = "Colwyn Bay"
$row172[0,1] = "Aberystwyth"
$row172[0,2] = 1
$row172[0,3] = 2
$row172[0,4] = "A"
$row172[0,5] = 2.1
$row172[0,6] = 3.2
$row172[0,7] = 3.1
$row172[0,8] = 2.3
$row172[0,9] = 3.25
$row172[0,10] = 2.75
$row172[0,11] = -0.25
$row172[0,12] = 2.025
$row172[0,13] = 1.775
$row172[0,14] = 2.5
$row172[0,15] = 2
$row172[0,16] = 1.8
$row172[0,17] = -1
$row172[0,18] = -1
$row172[0,19] = 1.75
$row172[0,20] = -1
$row172[0,21] = 0.7749999999999999
$row172[0,22] = 1
$row172[0,23] = -1
$ws.Range("F172:AC172").Value = $row172

$ws.Range("B173").Value = 7721622
$row173 = New-Object "object[,]" 1,24
$row173[0,0] = "Pontypridd Town"
$row173[0,1] = "Haverfordwest County"
$row173[0,2] = 0
$row173[0,3] = 1
$row173[0,4] = "A"
$row173[0,5] = 2.375
$row173[0,6] = 3.25
$row173[0,7] = 2.6
$row173[0,8] = 3
$row173[0,9] = 3.2
$row173[0,10] = 2.25
$row173[0,11] = 0.25
$row173[0,12] = 1.85
$row173[0,13] = 1.95
$row173[0,14] = 2.25
$row173[0,15] = 1.95
$row173[0,16] = 1.85
$row173[0,17] = -1
$row173[0,18] = -1
$row173[0,19] = 1.25
$row173[0,20] = -1
$row173[0,21] = 0.95
$row173[0,22] = -1
$row173[0,23] = 0.8500000000000001
$ws.Range("F173:AC173").Value = $row173

